# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to Famfrit_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 500
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("N16").Value = -960
$ws.Range("H17").Value = 1484058.1
$ws.Range("J17").Value = 1484058.1
$ws.Range("L17").Value = 4452174.300000001
$ws.Range("N17").Value = -4452510.300000001
$ws.Range("H48").Value = 3119
$ws.Range("J48").Value = 3119
$ws.Range("L48").Value = 9357
$ws.Range("N48").Value = -9941
$ws.Range("H53").Value = 988.7646999999999
$ws.Range("I53").Value = 774.1539
$ws.Range("K53").Value = 774.1539
$ws.Range("M53").Value = -137.1539
$ws.Range("H56").Value = 3119
$ws.Range("J56").Value = 3119
$ws.Range("L56").Value = 9357
$ws.Range("N56").Value = -10425
$ws.Range("H137").Value = 3141.6072
$ws.Range("I137").Value = 2958.2632
$ws.Range("J137").Value = 3528.6667
$ws.Range("K137").Value = 8874.7896
$ws.Range("L137").Value = 10586.0001
$ws.Range("M137").Value = -6324.7896
$ws.Range("N137").Value = -15686.0001
$ws.Range("H141").Value = 1976.5264
$ws.Range("I141").Value = 1884
$ws.Range("J141").Value = 2470
$ws.Range("K141").Value = 5652
$ws.Range("L141").Value = 7410
$ws.Range("M141").Value = -472
$ws.Range("N141").Value = -17770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 538.5714
$ws.Range("I2").Value = 484.5
$ws.Range("K2").Value = 484.5
$ws.Range("M2").Value = -371.5
$ws.Range("H74").Value = 58825148
$ws.Range("I74").Value = 62501564
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 62501564
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -62500690
$ws.Range("N74").Value = -4248
$ws.Range("H77").Value = 58825148
$ws.Range("I77").Value = 62501564
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 312507820
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -312503452
$ws.Range("N77").Value = -21236
$ws.Range("H116").Value = 538.5714
$ws.Range("I116").Value = 484.5
$ws.Range("K116").Value = 484.5
$ws.Range("M116").Value = 1809.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 538.5714
$ws.Range("I3").Value = 484.5
$ws.Range("K3").Value = 484.5
$ws.Range("M3").Value = -370.5
$ws.Range("H99").Value = 6252.5
$ws.Range("I99").Value = 2999
$ws.Range("K99").Value = 2999
$ws.Range("M99").Value = -1501
$ws.Range("H134").Value = 2559.6099
$ws.Range("I134").Value = 2484.8438
$ws.Range("K134").Value = 7454.5314
$ws.Range("M134").Value = -4919.5314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1446.6
$ws.Range("I16").Value = 1252.75
$ws.Range("J16").Value = 2222
$ws.Range("K16").Value = 1252.75
$ws.Range("L16").Value = 2222
$ws.Range("M16").Value = -965.75
$ws.Range("N16").Value = -2796
$ws.Range("H23").Value = 8500
$ws.Range("I23").Value = 8500
$ws.Range("K23").Value = 8500
$ws.Range("M23").Value = -8260
$ws.Range("H27").Value = 8500
$ws.Range("I27").Value = 8500
$ws.Range("K27").Value = 8500
$ws.Range("M27").Value = -8308
$ws.Range("H31").Value = 5110.125
$ws.Range("I31").Value = 3358.0625
$ws.Range("J31").Value = 8614.25
$ws.Range("K31").Value = 3358.0625
$ws.Range("L31").Value = 8614.25
$ws.Range("M31").Value = -3063.0625
$ws.Range("N31").Value = -9204.25
$ws.Range("H34").Value = 5110.125
$ws.Range("I34").Value = 3358.0625
$ws.Range("J34").Value = 8614.25
$ws.Range("K34").Value = 3358.0625
$ws.Range("L34").Value = 8614.25
$ws.Range("M34").Value = -3156.0625
$ws.Range("N34").Value = -9018.25
$ws.Range("H113").Value = 1446.6
$ws.Range("I113").Value = 1252.75
$ws.Range("J113").Value = 2222
$ws.Range("K113").Value = 1252.75
$ws.Range("L113").Value = 2222
$ws.Range("M113").Value = 917.25
$ws.Range("N113").Value = -6562

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 95.22221999999999
$ws.Range("I38").Value = 27.75
$ws.Range("J38").Value = 149.2
$ws.Range("K38").Value = 83.25
$ws.Range("L38").Value = 447.6
$ws.Range("M38").Value = 263.75
$ws.Range("N38").Value = -1141.6
$ws.Range("H56").Value = 21736.818
$ws.Range("I56").Value = 21736.818
$ws.Range("K56").Value = 21736.818
$ws.Range("M56").Value = -21206.818
$ws.Range("H104").Value = 7000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 7000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 21000
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -26242
$ws.Range("H121").Value = 600
$ws.Range("H131").Value = 1338.8096
$ws.Range("J131").Value = 1519.4546
$ws.Range("L131").Value = 4558.3638
$ws.Range("N131").Value = -14638.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H132").Value = 4355.778
$ws.Range("I132").Value = 2981.5
$ws.Range("K132").Value = 8944.5
$ws.Range("M132").Value = -6414.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2541.7334
$ws.Range("I40").Value = 2630.6924
$ws.Range("J40").Value = 1963.5
$ws.Range("K40").Value = 2630.6924
$ws.Range("L40").Value = 1963.5
$ws.Range("M40").Value = -2494.6924
$ws.Range("N40").Value = -2235.5
$ws.Range("H136").Value = 3349.75
$ws.Range("I136").Value = 2333
$ws.Range("J136").Value = 6400
$ws.Range("K136").Value = 6999
$ws.Range("L136").Value = 19200
$ws.Range("M136").Value = -4449
$ws.Range("N136").Value = -24300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1312.625
$ws.Range("I81").Value = 1099.8334
$ws.Range("K81").Value = 2199.6668
$ws.Range("M81").Value = -1138.6668
$ws.Range("H84").Value = 1312.625
$ws.Range("I84").Value = 1099.8334
$ws.Range("K84").Value = 10998.334
$ws.Range("M84").Value = -5694.333999999999
$ws.Range("H107").Value = 649.5
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840
$ws.Range("H136").Value = 2930.261
$ws.Range("I136").Value = 1022
$ws.Range("K136").Value = 3066
$ws.Range("M136").Value = -516
